$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.763.01"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "3.102.09"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'576.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'172.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.097.83"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'6.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "'0.478"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "'37.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "3.619.34"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "66.782.95"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "3.104.59"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'16.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").Value = "'476.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("D22").Value = "'0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "'7.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.82%  "
$ws.Range("D24").Value = "'13.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.20%  "
$ws.Range("D25").Value = "'83.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'9.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").Value = "'2.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").Value = "0.0₃0929"
$ws.Range("E34").Value = "  -8.96%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'5.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'0.975"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").Value = "'47.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").Value = "'2.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").Value = "'49.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "2.792.14"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "'0.0353"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").Value = "'376.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "'2.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.34%  "
$ws.Range("D48").Value = "'135.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'2.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.02%  "
